# Apply the "update data prep + model run" change set:
#  - Clear Input2 (C5) and Output2 (E5) on the Connections sheet, row 5
#  - Change Connection_type (F5) from "connection_type_lossless_bidirectional"
#    to the new value "connection_type_normal"
#  - Clear Cap_Input2_existing (I5) and Cap_Output2_existing (M5)
#  - Switch the active sheet/selection: Units becomes the active tab with
#    selection Q37; Connections loses the active tab flag and its selection
#    moves to I6

$wb = $excel.ActiveWorkbook

$wsUnits = $wb.Worksheets.Item("Units")
$wsConn  = $wb.Worksheets.Item("Connections")

# --- Data edits on the Connections sheet, row 5 ---
$wsConn.Range("C5").ClearContents()
$wsConn.Range("E5").ClearContents()
$wsConn.Range("F5").Value = "connection_type_normal"
$wsConn.Range("I5").ClearContents()
$wsConn.Range("M5").ClearContents()

# --- Sheet view / selection updates ---
# Make Connections active first and set its new selection, then switch
# to Units last so that it ends up being the active/selected tab.
$wsConn.Activate()
$wsConn.Range("I6").Select()

$wsUnits.Activate()
$wsUnits.Range("Q37").Select()
